$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.002.14'
$ws.Range('E2').Value = '  -1.07%  '
$ws.Range('D3').Value = '1.895.84'
$ws.Range('E3').Value = '  -1.88%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7377'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.73'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.86%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.9984'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3103'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -2.55%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '26.46'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -4.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.06911'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -1.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.7729'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -0.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07945'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.56%  '
$ws.Range('D13').Value = '1.902.90'
$ws.Range('E13').Value = '  -1.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.232'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -2.29%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '91.57'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.91%  '
$ws.Range('D16').Value = '30.016.88'
$ws.Range('E16').Value = '  -1.00%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.15'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.84%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.809'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +1.51%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '240.05'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -5.02%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007788'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.74%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.9986'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.17%  '
$ws.Range('D22').Value = '2.139.33'
$ws.Range('E22').Value = '  -1.92%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9981'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.910'
$ws.Range('D24').ClearFormats()
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.314'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '166.83'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.81'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.59%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.1279'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -3.65%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.034'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -8.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.355'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.48%  '
$ws.Range('E31').Value = '  +1.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.297'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.75%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.060'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.35%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.05117'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.77%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.281'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7372'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -1.00%  '
$ws.Range('E37').Value = '  -2.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01929'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.792'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +0.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.313'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '74.40'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.41%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.4459'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.940'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.25%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.9985'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8363'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '7.665'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +2.73%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '101.21'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '9.848'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('D49').Value = '2.045.56'
$ws.Range('E49').Value = '  -1.89%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '36.55'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '933.80'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -5.12%  '
